$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new sub-bullet paragraph right after the paragraph that
#    ends with "for your application" (and before the "If NuGet error..."
#    paragraph), containing the new firewall note.
# ------------------------------------------------------------------

$anchor = $d.Content
$found = $anchor.Find.Execute("If NuGet error", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'If NuGet error' paragraph to anchor the new bullet."
}

$anchor.Collapse(1)
$anchor.InsertParagraphBefore()

# The freshly inserted (still empty) paragraph is now the paragraph that
# immediately precedes the "If NuGet error..." paragraph.
$newPara = $anchor.Paragraphs.Item(1)
$newPara.Range.Text = "You may need to open port 9000 in the windows firewall as well"

# Match the indentation level used by the other sub-bullets (ilvl=1 in the
# OOXML, which is ListLevelNumber=2 in the Word object model).
$newPara.Range.ListFormat.ListLevelNumber = 2

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark so that it sits at the end of the new
#    paragraph's text (collapsed, right before the paragraph mark)
#    instead of at the end of the document.
# ------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$endRange = $newPara.Range.Duplicate
$endRange.Collapse(0)
$null = $endRange.MoveEnd(1, -1)

# Adding a bookmark to a genuinely zero-length range is unreliable, so we
# temporarily insert a marker character, bookmark the (non-empty) range
# that spans it, then delete the marker. The bookmark collapses back down
# to its correct zero-width position and survives the deletion.
$endRange.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $endRange)

$markerRange = $d.Range($endRange.Start, $endRange.End)
$markerRange.Delete()
